# Refresh crypto price/volume figures to reflect the latest scrape run.
# Numeric-looking "Price" strings are forced back to Text format (matching the
# source data, which stores prices like "67.79" as plain text) before assignment,
# otherwise Excel would auto-convert them into real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.030.42'
$ws.Range('E2').Value = '  +1.84%  '
$ws.Range('D3').Value = '2.211.42'
$ws.Range('E3').Value = '  +1.31%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '251.84'
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.622'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '67.79'
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').Value = '  -0.36%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.620'
$ws.Range('E9').Value = '  +7.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.93'
$ws.Range('E10').Value = '  +3.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '59.56'
$ws.Range('E11').Value = '  +2.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0938'
$ws.Range('E12').Value = '  +0.36%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.04'
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('E14').Value = '  +0.03%  '
$ws.Range('D15').Value = '2.539.46'
$ws.Range('E15').Value = '  +1.14%  '
$ws.Range('E16').Value = '  +0.63%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.50'
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').Value = '2.204.34'
$ws.Range('E18').Value = '  +0.48%  '
$ws.Range('D19').Value = '41.837.49'
$ws.Range('E19').Value = '  +1.60%  '
$ws.Range('D20').Value = '0.0₃0960'
$ws.Range('E20').Value = '  +1.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.26'
$ws.Range('E21').Value = '  +0.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.13'
$ws.Range('E22').Value = '  -1.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '231.27'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.02'
$ws.Range('E24').Value = '  -2.80%  '
$ws.Range('E25').Value = '  +1.45%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.16'
$ws.Range('E27').Value = '  -5.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.41'
$ws.Range('E28').Value = '  -3.39%  '
$ws.Range('E29').Value = '  -1.29%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.12'
$ws.Range('E30').Value = '  -2.34%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '166.56'
$ws.Range('E31').Value = '  -2.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.41'
$ws.Range('E32').Value = '  -0.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.122'
$ws.Range('E33').Value = '  +3.70%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.85'
$ws.Range('E34').Value = '  +7.87%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0780'
$ws.Range('E35').Value = '  +7.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.122'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.60'
$ws.Range('E37').Value = '  +0.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.71'
$ws.Range('E38').Value = '  +0.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.07'
$ws.Range('E39').Value = '  +1.92%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0308'
$ws.Range('E40').Value = '  +3.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.23'
$ws.Range('E41').Value = '  +0.51%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.20'
$ws.Range('E42').Value = '  +7.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.63'
$ws.Range('E43').Value = '  -0.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.00'
$ws.Range('E44').Value = '  -1.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.36'
$ws.Range('E45').Value = '  -4.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.196'
$ws.Range('E46').Value = '  -2.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.55'
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('E48').Value = '  -2.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.999'
$ws.Range('E49').Value = '  -0.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.15'
$ws.Range('E50').Value = '  +1.04%  '
$ws.Range('E51').Value = '  +4.34%  '
